# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold prices formatted as plain text (e.g. "72.144.08").
# Force a text number format first so purely-numeric-looking prices (like "536.21")
# are not silently re-interpreted as numbers when the value is assigned.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.144.08"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "4.021.53"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "536.21"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "152.36"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("E7").Value = "  +13.69%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").Value = "0.0000325"
$ws.Range("E11").Value = "  -6.73%  "
$ws.Range("D12").Value = "47.56"
$ws.Range("E12").Value = "  +11.23%  "
$ws.Range("D13").Value = "10.72"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "4.664.90"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "4.011.85"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "14.09"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "20.54"
$ws.Range("E17").Value = "  -3.51%  "
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").Value = "71.972.59"
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").Value = "428.85"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").Value = "99.00"
$ws.Range("E22").Value = "  +9.73%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  +5.04%  "
$ws.Range("D25").Value = "14.45"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").Value = "11.08"
$ws.Range("E26").Value = "  -8.41%  "
$ws.Range("D27").Value = "10.83"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("D29").Value = "36.93"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +25.07%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "679.69"
$ws.Range("E33").Value = "  -2.23%  "
$ws.Range("D34").Value = "6.98"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("D35").Value = "66.46"
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("D36").Value = "42.81"
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("D37").Value = "0.424"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0821"
$ws.Range("E39").Value = "  -11.42%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "3.39"
$ws.Range("E40").Value = "  +7.50%  "
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  +5.99%  "
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").Value = "9.63"
$ws.Range("E47").Value = "  +4.61%  "
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -11.25%  "
$ws.Range("E49").Value = "  -7.85%  "
$ws.Range("D50").Value = "3.33"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "146.36"
$ws.Range("E51").Value = "  +1.87%  "
